$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I and J, reusing the existing header format
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-25: values for columns I (I0) and J (IF)
$data = @(
    @(7, 7),
    @(6, 6),
    @(3, 4),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(4, 5),
    @(7, 7),
    @(9, 9),
    @(4, 5),
    @(5, 6),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(12, 12),
    @(6, 6),
    @(1, 2),
    @(1, 1),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(8, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
